$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A182").Value = "SEP0228264"
$ws.Range("B182").Value = 5500143583
$ws.Range("C182").Value = 2625
$ws.Range("D182").Value = "MPB-KÖL"
$ws.Range("E182").Value = "yes"

$ws.Range("A183").Value = "SEP0228263"
$ws.Range("B183").Value = 5500143583
$ws.Range("C183").Value = 2625
$ws.Range("D183").Value = "MPB-KÖL"
$ws.Range("E183").Value = "yes"

$ws.Range("A184").Value = "SEP0279425"
$ws.Range("B184").Value = 5500143583
$ws.Range("C184").Value = 2625
$ws.Range("D184").Value = "MPB-KÖL"
$ws.Range("E184").Value = "yes"

$ws.Range("A185").Value = "SEP0279423"
$ws.Range("B185").Value = 5500143583
$ws.Range("C185").Value = 2625
$ws.Range("D185").Value = "MPB-KÖL"
$ws.Range("E185").Value = "yes"

$ws.Range("A186").Value = "SEP0228262"
$ws.Range("B186").Value = 5500143583
$ws.Range("C186").Value = 2625
$ws.Range("D186").Value = "MPB-KÖL"
$ws.Range("E186").Value = "yes"

$ws.Range("A187").Value = "SEP0279421"
$ws.Range("B187").Value = 5500143583
$ws.Range("C187").Value = 2625
$ws.Range("D187").Value = "MPB-KÖL"
$ws.Range("E187").Value = "yes"

$ws.Range("A188").Value = "SEP0279422"
$ws.Range("B188").Value = 5500143583
$ws.Range("C188").Value = 2625
$ws.Range("D188").Value = "MPB-KÖL"
$ws.Range("E188").Value = "yes"

$ws.Range("A189").Value = "SEP0279424"
$ws.Range("B189").Value = 5500143583
$ws.Range("C189").Value = 2625
$ws.Range("D189").Value = "MPB-KÖL"
$ws.Range("E189").Value = "yes"

$ws.Range("A190").Value = "SEP0279420"
$ws.Range("B190").Value = 5500143583
$ws.Range("C190").Value = 2625
$ws.Range("D190").Value = "MPB-KÖL"
$ws.Range("E190").Value = "yes"

$ws.Range("A191").Value = "SEP0228265"
$ws.Range("B191").Value = 5500143583
$ws.Range("C191").Value = 2625
$ws.Range("D191").Value = "MPB-KÖL"
$ws.Range("E191").Value = "yes"

$ws.Range("A192").Value = "SEP0228266"
$ws.Range("B192").Value = 5500143583
$ws.Range("C192").Value = 2625
$ws.Range("D192").Value = "MPB-KÖL"
$ws.Range("E192").Value = "yes"

$ws.Range("A193").Value = "SEP0228267"
$ws.Range("B193").Value = 5500143583
$ws.Range("C193").Value = 2625
$ws.Range("D193").Value = "MPB-KÖL"
$ws.Range("E193").Value = "yes"

$ws.Range("A181:D193").Select()
